$wb = $excel.ActiveWorkbook

# ALC row 3
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 455223
$ws.Range("J3").Value = 455223
$ws.Range("L3").Value = 455223
$ws.Range("N3").Value = -455451

# ALC row 63
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# ALC row 66
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# ALC row 68
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 69999
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

# ALC row 71
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H71").Value = 69999
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6982.6665
$ws.Range("I76").Value = 6982.6665
$ws.Range("K76").Value = 6982.6665
$ws.Range("M76").Value = -6667.6665

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 6982.6665
$ws.Range("I79").Value = 6982.6665
$ws.Range("K79").Value = 6982.6665
$ws.Range("M79").Value = -5890.6665

# ALC row 102
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H102").Value = 455223
$ws.Range("J102").Value = 455223
$ws.Range("L102").Value = 455223
$ws.Range("N102").Value = -461713

# ALC row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 89899
$ws.Range("J136").Value = 89899
$ws.Range("L136").Value = 89899
$ws.Range("N136").Value = -100099

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5099.1875
$ws.Range("I138").Value = 2110.9092
$ws.Range("J138").Value = 6664.476
$ws.Range("K138").Value = 6332.7276
$ws.Range("L138").Value = 19993.428
$ws.Range("M138").Value = -1192.7276
$ws.Range("N138").Value = -30273.428

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4025.889
$ws.Range("I45").Value = 4122.3335
$ws.Range("K45").Value = 4122.3335
$ws.Range("M45").Value = -3745.3335

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 83429630
$ws.Range("I74").Value = 125141890
$ws.Range("K74").Value = 125141890
$ws.Range("M74").Value = -125141016

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 83429630
$ws.Range("I77").Value = 125141890
$ws.Range("K77").Value = 625709450
$ws.Range("M77").Value = -625705082

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 33718.332
$ws.Range("J110").Value = 9990.4
$ws.Range("L110").Value = 9990.4
$ws.Range("N110").Value = -14080.4

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 83805.78
$ws.Range("J133").Value = 83805.78
$ws.Range("L133").Value = 83805.78
$ws.Range("N133").Value = -88865.78

# ARM row 134
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 79997.336
$ws.Range("J134").Value = 79997.336
$ws.Range("L134").Value = 79997.336
$ws.Range("N134").Value = -90137.336

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 268.33334
$ws.Range("I22").Value = 262
$ws.Range("K22").Value = 262
$ws.Range("M22").Value = -89

# BSM row 59
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 71323.336
$ws.Range("I59").Value = 51985
$ws.Range("J59").Value = 110000
$ws.Range("K59").Value = 51985
$ws.Range("L59").Value = 110000
$ws.Range("M59").Value = -51138
$ws.Range("N59").Value = -111694

# BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 110000
$ws.Range("J132").Value = 110000
$ws.Range("L132").Value = 110000
$ws.Range("N132").Value = -120120

# BSM row 133
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 100780
$ws.Range("J133").Value = 100780
$ws.Range("L133").Value = 100780
$ws.Range("N133").Value = -110900

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4231.3477
$ws.Range("I134").Value = 4126.1
$ws.Range("K134").Value = 12378.3
$ws.Range("M134").Value = -9843.300000000001

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1422.5
$ws.Range("I16").Value = 1422.5
$ws.Range("K16").Value = 1422.5
$ws.Range("M16").Value = -1135.5

# CRP row 48
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 37737
$ws.Range("J48").Value = 37737
$ws.Range("L48").Value = 37737
$ws.Range("N48").Value = -38689

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3956.6667
$ws.Range("I62").Value = 3820
$ws.Range("K62").Value = 3820
$ws.Range("M62").Value = -3196

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3956.6667
$ws.Range("I65").Value = 3820
$ws.Range("K65").Value = 19100
$ws.Range("M65").Value = -15980

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1422.5
$ws.Range("I113").Value = 1422.5
$ws.Range("K113").Value = 1422.5
$ws.Range("M113").Value = 747.5

# CRP row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 101412.35
$ws.Range("J140").Value = 101412.35
$ws.Range("L140").Value = 101412.35
$ws.Range("N140").Value = -111772.35

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 96533.17999999999
$ws.Range("J141").Value = 103286.5
$ws.Range("L141").Value = 103286.5
$ws.Range("N141").Value = -113646.5

# CUL row 18
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 2642.111
$ws.Range("I18").Value = 539.8570999999999
$ws.Range("K18").Value = 1619.5713
$ws.Range("M18").Value = -1450.5713

# CUL row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 2479.3
$ws.Range("I99").Value = 1465.5
$ws.Range("K99").Value = 4396.5
$ws.Range("M99").Value = -2150.5

# CUL row 108
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 117.333336
$ws.Range("I108").Value = 117.333336
$ws.Range("K108").Value = 352.000008
$ws.Range("M108").Value = 2527.999992

# CUL row 119
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 3149.6
$ws.Range("I119").Value = 2187.25
$ws.Range("K119").Value = 6561.75
$ws.Range("M119").Value = -1723.75

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2285.4666
$ws.Range("I139").Value = 2098.6667
$ws.Range("K139").Value = 6296.000100000001
$ws.Range("M139").Value = -1156.000100000001

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1576.8235
$ws.Range("J22").Value = 1246.091
$ws.Range("L22").Value = 1246.091
$ws.Range("N22").Value = -1836.091

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1576.8235
$ws.Range("J27").Value = 1246.091
$ws.Range("L27").Value = 1246.091
$ws.Range("N27").Value = -1460.091

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1469.625
$ws.Range("I46").Value = 981.9231
$ws.Range("J46").Value = 3583
$ws.Range("K46").Value = 981.9231
$ws.Range("L46").Value = 3583
$ws.Range("M46").Value = -793.9231
$ws.Range("N46").Value = -3959

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 649.6539
$ws.Range("I55").Value = 421.5
$ws.Range("K55").Value = 421.5
$ws.Range("M55").Value = -248.5

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3423.8
$ws.Range("I61").Value = 2679.0908
$ws.Range("K61").Value = 2679.0908
$ws.Range("M61").Value = -2477.0908

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3423.8
$ws.Range("I113").Value = 2679.0908
$ws.Range("K113").Value = 2679.0908
$ws.Range("M113").Value = -509.0907999999999

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2196718.5
$ws.Range("I122").Value = 3554.1282
$ws.Range("J122").Value = 6948574.5
$ws.Range("K122").Value = 10662.3846
$ws.Range("L122").Value = 20845723.5
$ws.Range("M122").Value = -8212.384600000001
$ws.Range("N122").Value = -20850623.5

# LTW row 131
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 74895.2
$ws.Range("I131").Value = 52648
$ws.Range("J131").Value = 89726.664
$ws.Range("K131").Value = 52648
$ws.Range("L131").Value = 89726.664
$ws.Range("M131").Value = -47608
$ws.Range("N131").Value = -99806.664

# LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 87998
$ws.Range("I140").Value = 87997
$ws.Range("J140").Value = 87998.5
$ws.Range("K140").Value = 87997
$ws.Range("L140").Value = 87998.5
$ws.Range("M140").Value = -82817
$ws.Range("N140").Value = -98358.5

# WVR row 47
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 37495
$ws.Range("J47").Value = 37495
$ws.Range("L47").Value = 37495
$ws.Range("N47").Value = -38639

# WVR row 70
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 21776
$ws.Range("J70").Value = 21776
$ws.Range("L70").Value = 21776
$ws.Range("N70").Value = -22406

# WVR row 73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 21776
$ws.Range("J73").Value = 21776
$ws.Range("L73").Value = 21776
$ws.Range("N73").Value = -23960

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 72357.09
$ws.Range("J123").Value = 70542.8
$ws.Range("L123").Value = 70542.8
$ws.Range("N123").Value = -80342.8
